# "User percentage breakdowns can now be specified"
#
# Adds a new "specified_breakdowns" column (I) to the Trend_instructions
# sheet (the active sheet), flips the sample row's use_as_trend flag
# (column G) from "T" to "F", and sets the new column's sample value to "T".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in I1.
$ws.Range("I1").Value = "specified_breakdowns"

# use_as_trend (G2) is no longer "T" now that a breakdown is specified.
$ws.Range("G2").Value = "F"

# Sample value for the new specified_breakdowns column.
$ws.Range("I2").Value = "T"

# Match the author's final selection, resting on the new header cell.
$ws.Range("I1").Select()
